$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '30.718.78'
$ws.Range("E2").Value = '  +1.70%  '
$ws.Range("D3").Value = "'" + '1.899.47'
$ws.Range("E3").Value = '  +2.62%  '
$ws.Range("D4").Value = "'" + '1.0000'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = "'" + '239.03'
$ws.Range("E5").Value = '  +1.10%  '
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").Value = "'" + '0.4817'
$ws.Range("E7").Value = '  +0.98%  '
$ws.Range("E8").Value = '  +1.01%  '
$ws.Range("D9").Value = "'" + '0.06563'
$ws.Range("E9").Value = '  +1.38%  '
$ws.Range("D10").Value = "'" + '1.966.09'
$ws.Range("E10").Value = '  +6.13%  '
$ws.Range("D11").Value = "'" + '0.07459'
$ws.Range("E11").Value = '  +2.11%  '
$ws.Range("D12").Value = "'" + '16.74'
$ws.Range("E12").Value = '  +2.33%  '
$ws.Range("D13").Value = "'" + '5.125'
$ws.Range("E13").Value = '  -0.21%  '
$ws.Range("D14").Value = "'" + '88.16'
$ws.Range("E14").Value = '  +1.02%  '
$ws.Range("D15").Value = "'" + '0.6685'
$ws.Range("E15").Value = '  +3.57%  '
$ws.Range("D16").Value = "'" + '30.691.24'
$ws.Range("E16").Value = '  +1.82%  '
$ws.Range("D17").Value = "'" + '13.34'
$ws.Range("E17").Value = '  +0.75%  '
$ws.Range("D18").Value = "'" + '0.9997'
$ws.Range("E18").Value = '  +0.01%  '
$ws.Range("D19").Value = "'" + '0.000007628'
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").Value = "'" + '231.85'
$ws.Range("E20").Value = '  +3.20%  '
$ws.Range("D21").Value = "'" + '2.159.66'
$ws.Range("E21").Value = '  +2.86%  '
$ws.Range("D22").Value = "'" + '5.298'
$ws.Range("E22").Value = '  +0.48%  '
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").Value = "'" + '6.236'
$ws.Range("E24").Value = '  +2.21%  '
$ws.Range("D25").Value = "'" + '170.16'
$ws.Range("E25").Value = '  +3.79%  '
$ws.Range("D26").Value = "'" + '9.329'
$ws.Range("E26").Value = '  +1.32%  '
$ws.Range("E27").Value = '  +1.44%  '
$ws.Range("E28").Value = '  +2.92%  '
$ws.Range("D29").Value = "'" + '1.403'
$ws.Range("E29").Value = '  -1.64%  '
$ws.Range("D30").Value = "'" + '0.1019'
$ws.Range("E30").Value = '  +10.80%  '
$ws.Range("D31").Value = "'" + '4.371'
$ws.Range("E31").Value = '  +2.92%  '
$ws.Range("D32").Value = "'" + '4.039'
$ws.Range("E32").Value = '  +1.93%  '
$ws.Range("D33").Value = "'" + '0.05109'
$ws.Range("E33").Value = '  +1.76%  '
$ws.Range("D34").Value = "'" + '1.220'
$ws.Range("E34").Value = '  +7.18%  '
$ws.Range("D35").Value = "'" + '0.7572'
$ws.Range("E35").Value = '  +2.01%  '
$ws.Range("D36").Value = "'" + '2.713'
$ws.Range("E36").Value = '  +1.11%  '
$ws.Range("D37").Value = "'" + '0.01882'
$ws.Range("E37").Value = '  +3.04%  '
$ws.Range("E38").Value = '  +2.11%  '
$ws.Range("D39").Value = "'" + '0.9214'
$ws.Range("E39").Value = '  +2.13%  '
$ws.Range("E40").Value = '  +1.23%  '
$ws.Range("D41").Value = "'" + '107.16'
$ws.Range("E41").Value = '  +0.41%  '
$ws.Range("D42").Value = "'" + '0.4314'
$ws.Range("E42").Value = '  +1.46%  '
$ws.Range("D44").Value = "'" + '5.758'
$ws.Range("E44").Value = '  -2.87%  '
$ws.Range("D45").Value = "'" + '7.446'
$ws.Range("E45").Value = '  +0.18%  '
$ws.Range("D46").Value = "'" + '64.39'
$ws.Range("E46").Value = '  +0.75%  '
$ws.Range("E47").Value = '  -2.71%  '
$ws.Range("E48").Value = '  -4.37%  '
$ws.Range("E49").Value = '  +3.03%  '
$ws.Range("E50").Value = '  -1.23%  '
$ws.Range("D51").Value = "'" + '0.05676'
$ws.Range("E51").Value = '  +0.02%  '
